# Converts the generic column/row index schedule into one labeled with
# weekday names (row 1) and class start times (column A), and updates a
# few class entries that moved between time slots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: weekday headers ---
$ws.Range("B1").Value = "segunda"
$ws.Range("C1").Value = "terça"
$ws.Range("D1").Value = "quarta"
$ws.Range("E1").Value = "quinta"
$ws.Range("F1").Value = "sexta"

# --- Column A: time-of-day labels ---
$ws.Range("A2").Value  = "7:00"
$ws.Range("A3").Value  = "7:50"
$ws.Range("A4").Value  = "8:40"
$ws.Range("A5").Value  = "9:30"
$ws.Range("A6").Value  = "10:40"
$ws.Range("A7").Value  = "11:30"
$ws.Range("A8").Value  = "13:00"
$ws.Range("A9").Value  = "13:50"
$ws.Range("A10").Value = "14:40"
$ws.Range("A11").Value = "15:30"
$ws.Range("A12").Value = "16:40"
$ws.Range("A13").Value = "17:30"

# --- Class entries that shifted to the new, correct time slot ---
$ws.Range("C2").Value  = "Desenho Técnico"
$ws.Range("C4").Value  = "-"

$ws.Range("D3").Value  = "EAP"

$ws.Range("C7").Value  = "Circuitos Elétricos 2"
$ws.Range("C12").Value = "-"

$ws.Range("E8").Value  = "-"
$ws.Range("F8").Value  = "EAP"

$ws.Range("B9").Value  = "-"
